# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Mon Mar 20 05:45:21 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.551.59"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.762.79"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.69"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3836"
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3401"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.90"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.137"
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07382"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.35"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.334"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.764.09"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.038"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001072"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06659"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.05"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.28"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.365"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.553.13"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.99"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.56"
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.420"
$ws.Range("E27").Value = "  -4.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.408"
$ws.Range("E28").Value = "  -5.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.21"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.35"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.963.88"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.084"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.956"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08789"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.67"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02406"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6755"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.306"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06304"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2173"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.510"
$ws.Range("E41").Value = "  -8.78%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.239"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.216"
$ws.Range("E43").Value = "  -6.65%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.09"
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6229"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.824"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.44"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.074"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07366"
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.145"
$ws.Range("E51").Value = "  +2.42%  "
